# Rename the worksheet (sheet name: "data" -> "Sheet1")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Sheet1"

# Fix typo in LastName for Fahima: "Eldarrat" -> "Eldarat"
$ws.Range("B5").Value = "Eldarat"

# Add a new "Salary" column
$ws.Range("D1").Value = "Salary"
$ws.Range("D2").Value = 111000
$ws.Range("D3").Value = 118000
$ws.Range("D4").Value = 148000

# Update selection to D3 to match final state
$ws.Range("D3").Select()
